# Auto-generated Word COM-interop script to apply the cover-letter revision diff.
# Each top-level paragraph keeps its identity/position; only the paragraph's
# contents (runs / pPr) are swapped out via Range.InsertXML, which is the
# reliable way in this runtime to get exact run-level formatting (superscript,
# subscript, specific rFonts) that plain Find/Replace or Range.Text can't express.
$d = $word.ActiveDocument

function Set-ParaXml($para, $innerXml) {
    $payload = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p>" + $innerXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    [void]$para.Range.InsertXML($payload)
}

# Work from the last paragraph to the first so earlier edits do not shift the
# character offsets of paragraphs we still need to touch.

# --- Paragraph 5: data/code availability + conflicts-of-interest statement ---
$p5Xml = '<w:r><w:t>The data and all scripts for generating the figures are available from</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Zenodo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>https://doi.org/10.5281/zenodo.4959705</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t>. The authors have no conflicts of interest to declare.</w:t></w:r>'
Set-ParaXml $d.Paragraphs.Item(5) $p5Xml

# --- Paragraph 3: summary of revisions made for the resubmission ---
$p3Xml = '<w:r><w:t>We have substantially revised the manuscript from the original version. We feel that such substantial revisions were necessary to adequately respond to the comments and suggestions from the reviewers. Specifically, the reviewers recommended that we restructure the introduction and the methods section, as well as cut down the results and clarify the discussion, all of which we have done. The figures are mostly unchanged, but we did adjust the paramete</w:t></w:r><w:r><w:t>rization of the model in Fig. 6, as well as making some minor visual adjustments.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>The tables have been expanded in response to the reviewer comments, as we did not include all of the data from all three of our experiments initially. Finally</w:t></w:r><w:r><w:t xml:space="preserve">, in response to Reviewer 2’s comments, we have redone the statistical analysis. In addition to the updated analysis presented in the main text, we have provided an alternative analysis </w:t></w:r><w:r><w:t>in the supplemental materials in support of our main findings.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>'
Set-ParaXml $d.Paragraphs.Item(3) $p3Xml

# --- Paragraph 2: was a bare empty paragraph, now carries a tab stop ---
$p2Xml = '<w:pPr><w:tabs><w:tab w:val="left" w:pos="2508"/></w:tabs></w:pPr>'
Set-ParaXml $d.Paragraphs.Item(2) $p2Xml

# --- Paragraph 1: opening thank-you / abstract paragraph ---
$p1Xml = '<w:r><w:t xml:space="preserve">We thank you </w:t></w:r><w:r><w:t xml:space="preserve">again </w:t></w:r><w:r><w:t xml:space="preserve">for your time and your consideration of our work. </w:t></w:r><w:r><w:t>As a reminder, t</w:t></w:r><w:r><w:t>his research article outlines a new technique for constraining soil carbon models using the radiocarbon signature of heterotrophic respiration</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>∆</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>14</w:t></w:r><w:r><w:t>C-CO</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t xml:space="preserve">measured in laboratory incubations of archived soils. </w:t></w:r><w:r><w:t>Our major finding is that measuring ∆</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>14</w:t></w:r><w:r><w:t>C-CO</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> in incubations of archived soils is a promising </w:t></w:r><w:r><w:t>technique</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>for adding</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>powerful constraints on</w:t></w:r><w:r><w:t xml:space="preserve"> soil carbon models. However, we found that air-drying and rewetting mobilizes </w:t></w:r><w:r><w:t xml:space="preserve">slightly </w:t></w:r><w:r><w:t>older carbon as compared to what is</w:t></w:r><w:r><w:t xml:space="preserve"> respired in field-moist soils. The differences </w:t></w:r><w:r><w:t xml:space="preserve">in </w:t></w:r><w:r><w:t>∆</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>14</w:t></w:r><w:r><w:t>C-CO</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">were significant, </w:t></w:r><w:r><w:t>but the absolute impact on ∆</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>14</w:t></w:r><w:r><w:t>C-CO</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> is small, indicating limited bias for constraining models.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
Set-ParaXml $d.Paragraphs.Item(1) $p1Xml

Write-Output "done"
